$d = $word.ActiveDocument

# 1. Ativação date update
$d.Content.Find.Execute("Ativação: 01/01/2021", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2024", 2) | Out-Null

# 2. Objetivos paragraph: replace text, then add new italic (English) paragraph after it
$d.Content.Find.Execute("Apresentar os conceitos de Logística, Logística Reversa e Gestão da Cadeia de Suprimentos. Capacitar o aluno para aplicação de técnicas e métodos quantitativos para otimização dos problemas em Logística e Cadeias de Suprimentos.", $true, $false, $false, $false, $false, $true, 1, $false, "Fornecer conhecimentos que proporcionam uma visão holística a respeito da Cadeia de Suprimentos e da Logística, apresentando métodos e ferramentas para otimizar o desempenho das cadeias produtivas.", 2) | Out-Null

$objParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Fornecer conhecimentos que proporcionam")) {
        $objParaIndex = $i
        break
    }
}
$objPara = $d.Paragraphs($objParaIndex)
$newRange = $objPara.Range.InsertParagraphAfter()
$insertedPara = $d.Paragraphs($objParaIndex + 1)
$r = $insertedPara.Range
$r.InsertAfter("Provide knowledge that paves a holistic view of Logistics and Supply Chain, presenting methods and tools available to optimize the performance of production chains.")
$textLen = "Provide knowledge that paves a holistic view of Logistics and Supply Chain, presenting methods and tools available to optimize the performance of production chains.".Length
$txtRange = $d.Range($r.Start, $r.Start + $textLen)
$txtRange.Font.Italic = 1

# 3. Docente responsável
$d.Content.Find.Execute("5840560 - Marco Antonio Carvalho Pereira", $true, $false, $false, $false, $false, $true, 1, $false, "3295113 - José Eduardo Holler Branco", 2) | Out-Null

# 4. Programa resumido: replace text, then add new italic (English) paragraph after it
$d.Content.Find.Execute("1. Introdução: 2. Gestão estratégica3. Gestão dos relacionamentos4. Gestão global de suprimentos5. Avaliação de desempenho6. Mapeamento e análise de processos7. Gestão de demanda8. Gestão e coordenação de estoques9. Gestão da logística10. Logística reversa", $true, $false, $false, $false, $false, $true, 1, $false, "Gerenciamento da cadeia de suprimentos e da logística: planejamento, otimização e controle.", 2) | Out-Null

$progResIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Gerenciamento da cadeia de suprimentos e da logística")) {
        $progResIndex = $i
        break
    }
}
$progResPara = $d.Paragraphs($progResIndex)
$progResPara.Range.InsertParagraphAfter()
$insertedPara2 = $d.Paragraphs($progResIndex + 1)
$r2 = $insertedPara2.Range
$r2.InsertAfter("Supply chain and logistics management: planning, optimization and control.")
$textLen2 = "Supply chain and logistics management: planning, optimization and control.".Length
$txtRange2 = $d.Range($r2.Start, $r2.Start + $textLen2)
$txtRange2.Font.Italic = 1

# 5. Programa (long text): replace text, then add new italic (English) paragraph after it
$oldPrograma = "1. Introdução: A concorrência entre cadeias de suprimento. Definição operacional. A globalização e a gestão de cadeia de suprimentos. Governança das cadeias de suprimentos2. Gestão estratégica: Estratégia de cadeia de suprimentos. Produtos funcionais x produtos inovadores. Fluxos empurrados puxados e híbridos. Custo de transação e a decisão estratégica de comprar ou fazer. Padronização. Integração de parceiros da cadeia de suprimento no projeto de novos produtos e processos.3. Gestão dos relacionamentos: Confiança entre parceiros. Negociação. Gestão do relacionamento com clientes. Segmentação de produtos. Gestão do relacionamento com fornecedores4. Gestão global de suprimentos: Tipos de suplemento. Estrutura organizacional para suprimentos. O processo de suprimento. Coopetição. Ética e responsabilidade social na gestão global de suprimentos5. Avaliação de desempenho: O que é medição de desempenho? Porque medir desempenho. Características de uma boa medida de desempenho. Alinhamento de incentivos em cadeias globais de suprimento. Tipos de contrato de relacionamento6. Mapeamento e análise de processos: Principais processos na cadeia de suprimento. O modelo SCOR (Supply Chain Operations Reference). Análise e melhoramento de processos.7. Gestão de demanda: Ações sobre a demanda para redução de variabilidade. Causas da variabilidade da demanda. Previsão de demanda. Processo de previsão de vendas. Métodos usados em previsões. Método Delphi. Incerteza de previsão8. Gestão e coordenação de estoques: Definição de estoques. Causa do surgimento dos estoques. Tipos de estoque. VMI (vendor management inventory) - estoque gerenciado pelo distribuidor. VOI (vendor owner inventory) - consignação9. Gestão da logística: Centralização versus descentralização. Pontos de armazenagem/distribuição. Funções dos armazéns. Sistemas logísticos escalonados. Localização de unidades logísticas. Gestão de transportes na cadeia de suprimentos.10. Logística reversa: Conceito, importância, estrutura e tendências. Sustentabilidade. Ciclo fechado. Tipos de ciclo fechado. Motivação empresarial. Gerenciamento integrado de resíduos. Modelos de roteirização. Programação de frotas de veículos."
$newPrograma = "i) Introdução à Logística e Cadeia de Suprimentos; ii) Planejamento da cadeia de suprimentos; iii) Planejamento do transporte; iv) Custos logísticos; v) Tipos de cargas e sistemas de armazenamento; vi) Modelos de transporte; vii) Modelos de localização; viii) Planejamento do estoque; ix) Logística Reversa e Economia Circular; e x)  Controle da logística e cadeia de suprimentos."
$d.Content.Find.Execute($oldPrograma, $true, $false, $false, $false, $false, $true, 1, $false, $newPrograma, 2) | Out-Null

$programaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("i) Introdução à Logística e Cadeia de Suprimentos")) {
        $programaIndex = $i
        break
    }
}
$programaPara = $d.Paragraphs($programaIndex)
$programaPara.Range.InsertParagraphAfter()
$insertedPara3 = $d.Paragraphs($programaIndex + 1)
$r3 = $insertedPara3.Range
$newProgramaEn = "i) Introduction to Logistics and Supply Chain; ii) Supply chain planning; iii) Transport planning; iv) Logistic costs; v) Types of cargos and storage systems; vi) Stock planning; vii) Transport models; viii) Location models; ix) Reverse Logistics and Circular Economy; and x) Control of logistics and supply chain."
$r3.InsertAfter($newProgramaEn)
$txtRange3 = $d.Range($r3.Start, $r3.Start + $newProgramaEn.Length)
$txtRange3.Font.Italic = 1

# 6. Avaliação: Método / Critério / Norma de recuperação text replacements
$d.Content.Find.Execute("Aulas expositivas teóricas, aulas de exercícios.", $true, $false, $false, $false, $false, $true, 1, $false, "Provas, trabalhos em grupo, exercícios individuais e seminários.", 2) | Out-Null
$d.Content.Find.Execute("Média de Provas e trabalhos (MF).", $true, $false, $false, $false, $false, $true, 1, $false, "Média das atividades avaliativas.", 2) | Out-Null
$d.Content.Find.Execute("Prova de Recuperação (PR). A Nota final (NF) será a média aritmética entre MF e PR", $true, $false, $false, $false, $false, $true, 1, $false, "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação.", 2) | Out-Null

# 7. Bibliografia replacement
$oldBib = "CORRÊA, HENRIQUE LUIZ. Gestão de rede de suprimento: integrando cadeias de suprimento no mundo globalizado. Editora Atlas, 2009CORREA, HENRIQUE LUIZ. Administração de cadeias de suprimento e logística: o essencial. Editora Atlas 2014PIRES, SÉRGIO. Gestão da cadeia de suprimentos (Supply Chain Management): conceitos, estratégias, práticas e casos. Editora Atlas segunda edição. 2009IYER, ANANTH; SESHHADRI, SHIDHAR; VASHER, ROY. A gestão da cadeia de suprimentos da Toyota. Bookman. 2009MYERSON, PAUL. Lean Supply Chain and logistics management. McGrawHill. 2012"
$newBib = "BOWERSOX, D. J.; CLOSS, D. J.; COOPER, M. B.; BOWERSOX, J. C. Gestão Logística da Cadeia de Suprimentos. 4. ed. AMGH, 2013. 472 p.BARTHOLOMEU, D. B.; CAIXETA FILHO, J. V. Logística Ambiental de Resíduos Sólidos. São Paulo: Atlas, 2011, 249 p.CHOPRA, S.; MEINDL, P. Gestão da cadeia de suprimentos: estratégia, planejamento e Operações. 6. ed. Pearson, 2015. 544 p.CAIXETA FILHO, J. V.; MARTINS, R. S. (org.). Gestão Logística do Transporte de Cargas. São Paulo: Atlas, 2001. 296 p.CAIXETA FILHO, J. V.; GAMEIRO, A. H. (org.). Sistemas de Gerenciamento de Transporte: Modelagem Matemática. São Paulo: Atlas, 2001. 125 p.CAIXETA FILHO, J. V. Pesquisa Operacional: Técnicas de Otimização Aplicadas a Sistemas Agroindustriais. São Paulo: Atlas, 2001. 171 p.LEITE, P. R. Logística Reversa: Competividade e Sustentabilidade. 3. ed. São Paulo: Saraiva, 2017. 360 p."
$d.Content.Find.Execute($oldBib, $true, $false, $false, $false, $false, $true, 1, $false, $newBib, 2) | Out-Null

Write-Output "done"
